# Workbook / worksheet handles (already open).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: COMBUSTION_FIJA / DIESEL -> COMBUSTION_FIJA / GAS_NATURAL,
#     and the date moves from 2022-10-22 to 2022-04-01.
$ws.Range("B2").Value = "GAS_NATURAL"
$ws.Range("E2").Value = 44652

# --- Rows 3-6 keep their text, only the date moves from
#     2022-11-22 to 2022-05-01.
$ws.Range("E3").Value = 44682
$ws.Range("E4").Value = 44682
$ws.Range("E5").Value = 44682
$ws.Range("E6").Value = 44682

# --- The extra COMBUSTION_FIJA rows (KEROSENE, FUEL_OIL, NAFTA,
#     CARBON_DE_LEÑA, LEÑA) that used to live in rows 7-11 are gone now
#     - remove the rows (shifting nothing else up since they're the tail).
$ws.Range("A7:E11").Delete() | Out-Null

# --- Put the selection where the author left it on next open.
$ws.Range("H2").Select() | Out-Null

Write-Output "edit applied"
